$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Mangatepopo at d/s Intake"
$ws.Range("B2").Value = "Chlorophyll A"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = $False
$ws.Range("E2").Value = "ok"
$ws.Range("F2").Value = 0.922559112333533
$ws.Range("G2").Value = 0.0196078431372549
$ws.Range("H2").Value = 0.784313725490196
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5.5
$ws.Range("K2").Value = -0.714619565217391
$ws.Range("L2").Value = -1.34099268223529
$ws.Range("M2").Value = 0.06006286921909
$ws.Range("N2").Value = -12.9930830039526
$ws.Range("O2").Value = "RepSite"
$ws.Range("P2").Value = "Very likely improving"
$ws.Range("Q2").Value = 1820825.9
$ws.Range("R2").Value = 5674346.8
$ws.Range("S2").Value = "Ruapehu District"
$ws.Range("T2").Value = "Whanganui"
$ws.Range("U2").Value = "Upper Whanganui"
$ws.Range("V2").Value = "Whai_1"
$ws.Range("W2").Value = "mg/m2"

# Row 3
$ws.Range("A3").Value = "Mangatepopo at d/s Intake"
$ws.Range("B3").Value = "Dissolved Oxygen Concentration"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = $True
$ws.Range("E3").Value = "ok"
$ws.Range("F3").Value = 0.269627285543471
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.9019607843137259
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 10.78
$ws.Range("K3").Value = -0.0232643312101912
$ws.Range("L3").Value = -0.0994787740604119
$ws.Range("M3").Value = 0.0589729079448083
$ws.Range("N3").Value = -0.215810122543517
$ws.Range("O3").Value = "RepSite"
$ws.Range("P3").Value = "Unlikely increasing"
$ws.Range("Q3").Value = 1820825.9
$ws.Range("R3").Value = 5674346.8
$ws.Range("S3").Value = "Ruapehu District"
$ws.Range("T3").Value = "Whanganui"
$ws.Range("U3").Value = "Upper Whanganui"
$ws.Range("V3").Value = "Whai_1"
$ws.Range("W3").Value = "g/m3"

# Row 4
$ws.Range("A4").Value = "Mangatepopo at d/s Intake"
$ws.Range("B4").Value = "Dissolved Reactive Phosphorus"
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = $True
$ws.Range("E4").Value = "ok"
$ws.Range("F4").Value = 0.141565435331173
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.450980392156863
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.013
$ws.Range("K4").Value = 0.0009925271739130001
$ws.Range("L4").Value = -0.0002324660073018
$ws.Range("M4").Value = 0.001337912087912
$ws.Range("N4").Value = 7.63482441471572
$ws.Range("O4").Value = "RepSite"
$ws.Range("P4").Value = "Unlikely improving"
$ws.Range("Q4").Value = 1820825.9
$ws.Range("R4").Value = 5674346.8
$ws.Range("S4").Value = "Ruapehu District"
$ws.Range("T4").Value = "Whanganui"
$ws.Range("U4").Value = "Upper Whanganui"
$ws.Range("V4").Value = "Whai_1"
$ws.Range("W4").Value = "mg/L"

# Row 5
$ws.Range("A5").Value = "Mangatepopo at d/s Intake"
$ws.Range("B5").Value = "Ammoniacal Nitrogen (NH4)"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = $False
$ws.Range("E5").Value = "< 5 Non-censored values"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = 0.9411764705882349
$ws.Range("H5").Value = 0.07843137254901961
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = "RepSite"
$ws.Range("P5").Value = "Not Analysed improving"
$ws.Range("Q5").Value = 1820825.9
$ws.Range("R5").Value = 5674346.8
$ws.Range("S5").Value = "Ruapehu District"
$ws.Range("T5").Value = "Whanganui"
$ws.Range("U5").Value = "Upper Whanganui"
$ws.Range("V5").Value = "Whai_1"
$ws.Range("W5").Value = "mg/L"

# Row 6
$ws.Range("A6").Value = "Mangatepopo at d/s Intake"
$ws.Range("B6").Value = "Nitrite Nitrogen (NO2)"
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = $False
$ws.Range("E6").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F6").Value = 0.977895740179676
$ws.Range("G6").Value = 0.803921568627451
$ws.Range("H6").Value = 0.09803921568627449
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 0.001
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = "RepSite"
$ws.Range("P6").Value = "Extremely likely improving"
$ws.Range("Q6").Value = 1820825.9
$ws.Range("R6").Value = 5674346.8
$ws.Range("S6").Value = "Ruapehu District"
$ws.Range("T6").Value = "Whanganui"
$ws.Range("U6").Value = "Upper Whanganui"
$ws.Range("V6").Value = "Whai_1"
$ws.Range("W6").Value = "mg/L"

# Row 7
$ws.Range("A7").Value = "Mangatepopo at d/s Intake"
$ws.Range("B7").Value = "Nitrate Nitrogen (NO3)"
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = $False
$ws.Range("E7").Value = "ok"
$ws.Range("F7").Value = 0.012676606964227
$ws.Range("G7").Value = 0.215686274509804
$ws.Range("H7").Value = 0.392156862745098
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 0.008999999999999999
$ws.Range("K7").Value = 0.0017588282504012
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0.0028675415087161
$ws.Range("N7").Value = 19.5425361155698
$ws.Range("O7").Value = "RepSite"
$ws.Range("P7").Value = "Extremely unlikely improving"
$ws.Range("Q7").Value = 1820825.9
$ws.Range("R7").Value = 5674346.8
$ws.Range("S7").Value = "Ruapehu District"
$ws.Range("T7").Value = "Whanganui"
$ws.Range("U7").Value = "Upper Whanganui"
$ws.Range("V7").Value = "Whai_1"
$ws.Range("W7").Value = "mg/L"

# Row 8
$ws.Range("A8").Value = "Mangatepopo at d/s Intake"
$ws.Range("B8").Value = "pH"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = $False
$ws.Range("E8").Value = "ok"
$ws.Range("F8").Value = 0.000528279958919
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.745098039215686
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 7.58
$ws.Range("K8").Value = -0.0670867346938774
$ws.Range("L8").Value = -0.104357142857143
$ws.Range("M8").Value = -0.0401373626373627
$ws.Range("N8").Value = -0.885049270367776
$ws.Range("O8").Value = "RepSite"
$ws.Range("P8").Value = "Exceptionally unlikely increasing"
$ws.Range("Q8").Value = 1820825.9
$ws.Range("R8").Value = 5674346.8
$ws.Range("S8").Value = "Ruapehu District"
$ws.Range("T8").Value = "Whanganui"
$ws.Range("U8").Value = "Upper Whanganui"
$ws.Range("V8").Value = "Whai_1"
$ws.Range("W8").Value = ""

# Row 9
$ws.Range("A9").Value = "Mangatepopo at d/s Intake"
$ws.Range("B9").Value = "SIN (Soluble Inorganic nitrogen)"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = $True
$ws.Range("E9").Value = "ok"
$ws.Range("F9").Value = 0.0076611205275503
$ws.Range("G9").Value = 0.0392156862745098
$ws.Range("H9").Value = 0.568627450980392
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0.014
$ws.Range("K9").Value = 0.0019904632152588
$ws.Range("L9").Value = 0.0004982680100253
$ws.Range("M9").Value = 0.0035503732035046
$ws.Range("N9").Value = 14.2175943947061
$ws.Range("O9").Value = "RepSite"
$ws.Range("P9").Value = "Exceptionally unlikely improving"
$ws.Range("Q9").Value = 1820825.9
$ws.Range("R9").Value = 5674346.8
$ws.Range("S9").Value = "Ruapehu District"
$ws.Range("T9").Value = "Whanganui"
$ws.Range("U9").Value = "Upper Whanganui"
$ws.Range("V9").Value = "Whai_1"
$ws.Range("W9").Value = "g/m3"

# Row 10
$ws.Range("A10").Value = "Mangatepopo at d/s Intake"
$ws.Range("B10").Value = "Chlorophyll A"
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = $False
$ws.Range("E10").Value = "ok"
$ws.Range("F10").Value = 0.09299974879615561
$ws.Range("G10").Value = 0.0092592592592592
$ws.Range("H10").Value = 0.657407407407407
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 3.85
$ws.Range("K10").Value = 0.16641970722635
$ws.Range("L10").Value = -0.0253971753047344
$ws.Range("M10").Value = 0.451085384891948
$ws.Range("N10").Value = 4.322589798087
$ws.Range("O10").Value = "RepSite"
$ws.Range("P10").Value = "Very unlikely improving"
$ws.Range("Q10").Value = 1820825.9
$ws.Range("R10").Value = 5674346.8
$ws.Range("S10").Value = "Ruapehu District"
$ws.Range("T10").Value = "Whanganui"
$ws.Range("U10").Value = "Upper Whanganui"
$ws.Range("V10").Value = "Whai_1"
$ws.Range("W10").Value = "mg/m2"

# Row 11
$ws.Range("A11").Value = "Mangatepopo at d/s Intake"
$ws.Range("B11").Value = "Dissolved Oxygen Concentration"
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = $True
$ws.Range("E11").Value = "ok"
$ws.Range("F11").Value = 0.791172050256838
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0.811111111111111
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 10.765
$ws.Range("K11").Value = 0.019991789819376
$ws.Range("L11").Value = -0.0150089665633143
$ws.Range("M11").Value = 0.0413825110222806
$ws.Range("N11").Value = 0.185711006218077
$ws.Range("O11").Value = "RepSite"
$ws.Range("P11").Value = "Likely increasing"
$ws.Range("Q11").Value = 1820825.9
$ws.Range("R11").Value = 5674346.8
$ws.Range("S11").Value = "Ruapehu District"
$ws.Range("T11").Value = "Whanganui"
$ws.Range("U11").Value = "Upper Whanganui"
$ws.Range("V11").Value = "Whai_1"
$ws.Range("W11").Value = "g/m3"

# Row 12
$ws.Range("A12").Value = "Mangatepopo at d/s Intake"
$ws.Range("B12").Value = "Dissolved Reactive Phosphorus"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = $False
$ws.Range("E12").Value = "ok"
$ws.Range("F12").Value = 0.982466553558188
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.323232323232323
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0.015
$ws.Range("K12").Value = -0.000494248985115
$ws.Range("L12").Value = -0.0009022231370934
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = -3.29499323410014
$ws.Range("O12").Value = "RepSite"
$ws.Range("P12").Value = "Extremely likely improving"
$ws.Range("Q12").Value = 1820825.9
$ws.Range("R12").Value = 5674346.8
$ws.Range("S12").Value = "Ruapehu District"
$ws.Range("T12").Value = "Whanganui"
$ws.Range("U12").Value = "Upper Whanganui"
$ws.Range("V12").Value = "Whai_1"
$ws.Range("W12").Value = "mg/L"

# Row 13
$ws.Range("A13").Value = "Mangatepopo at d/s Intake"
$ws.Range("B13").Value = "Ammoniacal Nitrogen (NH4)"
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = $False
$ws.Range("E13").Value = "WARNING: Sen slope based on two censored values"
$ws.Range("F13").Value = 0.565881289282962
$ws.Range("G13").Value = 0.947368421052632
$ws.Range("H13").Value = 0.06315789473684209
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 0.005
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = "RepSite"
$ws.Range("P13").Value = "As likely as not improving"
$ws.Range("Q13").Value = 1820825.9
$ws.Range("R13").Value = 5674346.8
$ws.Range("S13").Value = "Ruapehu District"
$ws.Range("T13").Value = "Whanganui"
$ws.Range("U13").Value = "Upper Whanganui"
$ws.Range("V13").Value = "Whai_1"
$ws.Range("W13").Value = "mg/L"

# Row 14
$ws.Range("A14").Value = "Mangatepopo at d/s Intake"
$ws.Range("B14").Value = "Nitrite Nitrogen (NO2)"
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = $False
$ws.Range("E14").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F14").Value = 0.991101600566305
$ws.Range("G14").Value = 0.626262626262626
$ws.Range("H14").Value = 0.101010101010101
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 0.001
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = "RepSite"
$ws.Range("P14").Value = "Virtually certain improving"
$ws.Range("Q14").Value = 1820825.9
$ws.Range("R14").Value = 5674346.8
$ws.Range("S14").Value = "Ruapehu District"
$ws.Range("T14").Value = "Whanganui"
$ws.Range("U14").Value = "Upper Whanganui"
$ws.Range("V14").Value = "Whai_1"
$ws.Range("W14").Value = "mg/L"

# Row 15
$ws.Range("A15").Value = "Mangatepopo at d/s Intake"
$ws.Range("B15").Value = "Nitrate Nitrogen (NO3)"
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = $False
$ws.Range("E15").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F15").Value = 0.0977547528911442
$ws.Range("G15").Value = 0.303030303030303
$ws.Range("H15").Value = 0.373737373737374
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 0.008999999999999999
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0.0006012210924073
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = "RepSite"
$ws.Range("P15").Value = "Very unlikely improving"
$ws.Range("Q15").Value = 1820825.9
$ws.Range("R15").Value = 5674346.8
$ws.Range("S15").Value = "Ruapehu District"
$ws.Range("T15").Value = "Whanganui"
$ws.Range("U15").Value = "Upper Whanganui"
$ws.Range("V15").Value = "Whai_1"
$ws.Range("W15").Value = "mg/L"

# Row 16
$ws.Range("A16").Value = "Mangatepopo at d/s Intake"
$ws.Range("B16").Value = "pH"
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = $False
$ws.Range("E16").Value = "ok"
$ws.Range("F16").Value = 0.000430668566457
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0.526315789473684
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 7.64
$ws.Range("K16").Value = -0.0194134396355354
$ws.Range("L16").Value = -0.0285672126966327
$ws.Range("M16").Value = -0.0107077589235727
$ws.Range("N16").Value = -0.254102613030568
$ws.Range("O16").Value = "RepSite"
$ws.Range("P16").Value = "Exceptionally unlikely increasing"
$ws.Range("Q16").Value = 1820825.9
$ws.Range("R16").Value = 5674346.8
$ws.Range("S16").Value = "Ruapehu District"
$ws.Range("T16").Value = "Whanganui"
$ws.Range("U16").Value = "Upper Whanganui"
$ws.Range("V16").Value = "Whai_1"
$ws.Range("W16").Value = ""

# Row 17
$ws.Range("A17").Value = "Mangatepopo at d/s Intake"
$ws.Range("B17").Value = "SIN (Soluble Inorganic nitrogen)"
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = $True
$ws.Range("E17").Value = "ok"
$ws.Range("F17").Value = 0.140182542761315
$ws.Range("G17").Value = 0.101010101010101
$ws.Range("H17").Value = 0.525252525252525
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 0.014
$ws.Range("K17").Value = 0.000332347588717
$ws.Range("L17").Value = -0.0001429549902152
$ws.Range("M17").Value = 0.000736719146512
$ws.Range("N17").Value = 2.37391134797868
$ws.Range("O17").Value = "RepSite"
$ws.Range("P17").Value = "Unlikely improving"
$ws.Range("Q17").Value = 1820825.9
$ws.Range("R17").Value = 5674346.8
$ws.Range("S17").Value = "Ruapehu District"
$ws.Range("T17").Value = "Whanganui"
$ws.Range("U17").Value = "Upper Whanganui"
$ws.Range("V17").Value = "Whai_1"
$ws.Range("W17").Value = "g/m3"
